$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 78,3
$arr[0,0] = 3667.430990000001
$arr[0,1] = 1.402125977911055
$arr[0,2] = 37932994.34286775
$arr[1,0] = 2984.36849
$arr[1,1] = 11.43002836359665
$arr[1,2] = 77713315.70467912
$arr[2,0] = 692.1148999999999
$arr[2,1] = 10.80911589739844
$arr[2,2] = 656725.7406632416
$arr[3,0] = 5871.803
$arr[3,1] = 9.867354356683791
$arr[3,2] = 17141001.91086065
$arr[4,0] = 2022.20619
$arr[4,1] = 6.711032608058304
$arr[4,2] = 39055016.51785336
$arr[5,0] = 4629.59329
$arr[5,1] = 7.751758998725563
$arr[5,2] = 49484956.59412816
$arr[6,0] = 347.0675
$arr[6,1] = 9.733060268685222
$arr[6,2] = 51035962.90503629
$arr[7,0] = 6216.8504
$arr[7,1] = 10.07750778272748
$arr[7,2] = 88657102.460186
$arr[8,0] = 5610.70249
$arr[8,1] = 5.106852845288813
$arr[8,2] = 2221001.093334053
$arr[9,0] = 1041.09699
$arr[9,1] = 9.308235140517354
$arr[9,2] = 91666424.78696537
$arr[10,0] = 5910.8251
$arr[10,1] = 8.62972926069051
$arr[10,2] = 49656270.03224567
$arr[11,0] = 653.0928
$arr[11,1] = 7.873177503235638
$arr[11,2] = 57138227.63230186
$arr[12,0] = 5915.74379
$arr[12,1] = 2.5615228950046
$arr[12,2] = 63544596.3457576
$arr[13,0] = 736.0556899999999
$arr[13,1] = 13.92321351682767
$arr[13,2] = 69062905.55134416
$arr[14,0] = 4280.6112
$arr[14,1] = 4.306655411608517
$arr[14,2] = 52234534.49153807
$arr[15,0] = 2283.3067
$arr[15,1] = 6.38427272439003
$arr[15,2] = 80986035.04535043
$arr[16,0] = 4930.699900000001
$arr[16,1] = 12.4599950904958
$arr[16,2] = 13937808.8820586
$arr[17,0] = 1633.218
$arr[17,1] = 7.889897698536515
$arr[17,2] = 48799388.31891399
$arr[18,0] = 1348.10629
$arr[18,1] = 6.673630793578923
$arr[18,2] = 79433581.65045036
$arr[19,0] = 5303.69319
$arr[19,1] = 10.30259809689596
$arr[19,2] = 20496555.46071241
$arr[20,0] = 2635.3864
$arr[20,1] = 6.800201285164803
$arr[20,2] = 72999135.69416618
$arr[21,0] = 3928.531499999999
$arr[21,1] = 2.216097831260413
$arr[21,2] = 10419957.35912817
$arr[22,0] = 3972.472290000001
$arr[22,1] = 11.12857661768794
$arr[22,2] = 37208648.50963698
$arr[23,0] = 2679.32719
$arr[23,1] = 2.381598450243473
$arr[23,2] = 14092490.05015939
$arr[24,0] = 1677.15879
$arr[24,1] = 2.060205020941794
$arr[24,2] = 4907115.145509597
$arr[25,0] = 4974.64069
$arr[25,1] = 4.577458583749831
$arr[25,2] = 11377913.94786909
$arr[26,0] = 997.1562
$arr[26,1] = 14.36199652310461
$arr[26,2] = 6981399.813687894
$arr[27,0] = 5566.7617
$arr[27,1] = 12.81600432284176
$arr[27,2] = 76333689.18198627
$arr[28,0] = 4301.52479
$arr[28,1] = 10.89899302227423
$arr[28,2] = 93007715.94594233
$arr[29,0] = 2350.27469
$arr[29,1] = 4.260187649633735
$arr[29,2] = 3187405.073873233
$arr[30,0] = 6260.791190000001
$arr[30,1] = 7.993033485021442
$arr[30,2] = 64184014.28665267
$arr[31,0] = 391.00829
$arr[31,1] = 9.767404143698514
$arr[31,2] = 89083112.52863845
$arr[32,0] = 4585.6525
$arr[32,1] = 8.263915817253292
$arr[32,2] = 24004069.33576101
$arr[33,0] = 1978.2654
$arr[33,1] = 5.289341134950519
$arr[33,2] = 76942785.89646564
$arr[34,0] = 4607.550090000001
$arr[34,1] = 5.193924769293517
$arr[34,2] = 80814409.52896607
$arr[35,0] = 2044.24939
$arr[35,1] = 3.518531526438892
$arr[35,2] = 51721325.35526296
$arr[36,0] = 2327.24749
$arr[36,1] = 12.54478178638965
$arr[36,2] = 10759601.85667034
$arr[37,0] = 4324.55199
$arr[37,1] = 6.036846867762506
$arr[37,2] = 55993815.1162886
$arr[38,0] = 1671.2561
$arr[38,1] = 3.91916290903464
$arr[38,2] = 61335760.61696513
$arr[39,0] = 4892.6618
$arr[39,1] = 1.865453130099922
$arr[39,2] = 1438217.865468469
$arr[40,0] = 367.98109
$arr[40,1] = 8.283902868162841
$arr[40,2] = 38142768.8839566
$arr[41,0] = 6283.81839
$arr[41,1] = 5.686381690669805
$arr[41,2] = 91085292.70639411
$arr[42,0] = 697.03359
$arr[42,1] = 5.751517934724689
$arr[42,2] = 62225989.09788997
$arr[43,0] = 5954.765890000001
$arr[43,1] = 5.667865811381489
$arr[43,2] = 35967270.44749773
$arr[44,0] = 2306.3339
$arr[44,1] = 14.70432431669906
$arr[44,2] = 6199977.065476123
$arr[45,0] = 4257.584
$arr[45,1] = 6.674329224508256
$arr[45,2] = 65943896.64192218
$arr[46,0] = 2672.29489
$arr[46,1] = 14.03565506171435
$arr[46,2] = 79025006.68819761
$arr[47,0] = 3979.50459
$arr[47,1] = 12.50415098015219
$arr[47,2] = 75588452.68290257
$arr[48,0] = 2000.3086
$arr[48,1] = 7.814556941390038
$arr[48,2] = 6275993.930154946
$arr[49,0] = 4563.6093
$arr[49,1] = 13.63885116158053
$arr[49,2] = 24515605.72786955
$arr[50,0] = 3630.5225
$arr[50,1] = 10.70069581875578
$arr[50,2] = 77382643.0440764
$arr[51,0] = 2933.395399999999
$arr[51,1] = 6.206627360079437
$arr[51,2] = 35751964.84585525
$arr[52,0] = 1326.2087
$arr[52,1] = 3.22643582848832
$arr[52,2] = 22559245.18157728
$arr[53,0] = 5237.7092
$arr[53,1] = 14.4543857886456
$arr[53,2] = 16526074.07042058
$arr[54,0] = 5259.7524
$arr[54,1] = 12.51670508924872
$arr[54,2] = 99769304.85403165
$arr[55,0] = 1304.1655
$arr[55,1] = 2.896263035014272
$arr[55,2] = 50618585.26367228
$arr[56,0] = 4936.60259
$arr[56,1] = 14.76390654081479
$arr[56,2] = 37386445.26025513
$arr[57,0] = 1715.19689
$arr[57,1] = 13.83768974384293
$arr[57,2] = 98050597.0646441
$arr[58,0] = 3322.38359
$arr[58,1] = 10.64085420966148
$arr[58,2] = 86460740.57305464
$arr[59,0] = 3329.41589
$arr[59,1] = 2.630695700645447
$arr[59,2] = 95967952.19102409
$arr[60,0] = 6545.9029
$arr[60,1] = 9.514576449990273
$arr[60,2] = 47762697.00252451
$arr[61,0] = 5565.7777
$arr[61,1] = 9.464815934654325
$arr[61,2] = 95205019.72509316
$arr[62,0] = 998.1401999999999
$arr[62,1] = 3.792240146081895
$arr[62,2] = 25585956.61182189
$arr[63,0] = 3935.5638
$arr[63,1] = 13.48286427510902
$arr[63,2] = 24829822.58643489
$arr[64,0] = 2628.3541
$arr[64,1] = 9.718098845332861
$arr[64,2] = 4475884.721206967
$arr[65,0] = 6589.843690000001
$arr[65,1] = 5.352878645062447
$arr[65,2] = 7835591.810459737
$arr[66,0] = 2977.33619
$arr[66,1] = 1.171910731587559
$arr[66,2] = 12857010.16250066
$arr[67,0] = 3674.46329
$arr[67,1] = 1.612447910476476
$arr[67,2] = 10207254.45065415
$arr[68,0] = 3285.4751
$arr[68,1] = 9.149751935619861
$arr[68,2] = 74951090.73365107
$arr[69,0] = 3278.442799999999
$arr[69,1] = 14.99214413762093
$arr[69,2] = 63782077.68660271
$arr[70,0] = 5281.64999
$arr[70,1] = 5.767604018095881
$arr[70,2] = 96452203.14141828
$arr[71,0] = 1370.14949
$arr[71,1] = 1.045445160940289
$arr[71,2] = 95981702.70168921
$arr[72,0] = 1042.08099
$arr[72,1] = 11.3314286605455
$arr[72,2] = 2247760.809466708
$arr[73,0] = 5609.71849
$arr[73,1] = 4.739535705652088
$arr[73,2] = 14285239.48209127
$arr[74,0] = 6239.8776
$arr[74,1] = 2.827977422159165
$arr[74,2] = 87506651.43341525
$arr[75,0] = 324.0403
$arr[75,1] = 14.11988576920703
$arr[75,2] = 37267547.98032111
$arr[76,0] = 2940.4277
$arr[76,1] = 10.81613249657676
$arr[76,2] = 14411002.08812393
$arr[77,0] = 3623.490199999999
$arr[77,1] = 2.481001087930053
$arr[77,2] = 44898948.86223367

$ws.Range("A2:C79").Value = $arr
